$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" (D) and "Volume(1h)" (E) columns for rows 2-51 with latest crypto data.
# Some Price values look like plain numbers (e.g. "6.70", "1.00"); set the cell to Text
# format first so Excel keeps the exact text (incl. trailing zeros) instead of coercing it
# to a numeric value.

$ws.Range("D2").Value = "62.962.78"
$ws.Range("E2").Value = "  +4.22%  "
$ws.Range("D3").Value = "2.700.90"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.49"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.84"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").Value = "2.730.91"
$ws.Range("E9").Value = "  +4.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.113"
$ws.Range("E11").Value = "  +6.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  +4.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.157"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "3.183.72"
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.52"
$ws.Range("E15").Value = "  +8.45%  "
$ws.Range("D16").Value = "62.875.76"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000151"
$ws.Range("E17").Value = "  +6.85%  "
$ws.Range("D18").Value = "2.714.97"
$ws.Range("E18").Value = "  +3.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.94"
$ws.Range("E19").Value = "  +5.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.87"
$ws.Range("E20").Value = "  +4.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.15"
$ws.Range("E21").Value = "  +4.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.99"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.531"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.33"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +3.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.60"
$ws.Range("E27").Value = "  +7.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.03"
$ws.Range("E29").Value = "  +6.58%  "
$ws.Range("D30").Value = "0.0₃0851"
$ws.Range("E30").Value = "  +6.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.10"
$ws.Range("E31").Value = "  +11.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "170.07"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.54"
$ws.Range("E34").Value = "  +5.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.17"
$ws.Range("E35").Value = "  +18.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.75"
$ws.Range("E36").Value = "  +10.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.42"
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  +10.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.01"
$ws.Range("E39").Value = "  +20.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "350.79"
$ws.Range("E40").Value = "  +12.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.26"
$ws.Range("E41").Value = "  +9.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.17"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.68"
$ws.Range("E43").Value = "  +13.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.56"
$ws.Range("E44").Value = "  +8.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0594"
$ws.Range("E45").Value = "  +7.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.58"
$ws.Range("E46").Value = "  +8.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0260"
$ws.Range("E47").Value = "  +7.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.639"
$ws.Range("E48").Value = "  +5.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.08"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.996"
$ws.Range("E51").Value = "  -0.27%  "
